$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: capture the existing Circular Link URLs (column F) for rows 2..99 ---
# These rows shift down by one (to rows 3..100) once we insert the new row, so
# grab the current text now while the row numbers still match the original file.
$lastRow = 99
$urls = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $urls[$r] = $ws.Cells.Item($r, 6).Value()
}

# --- Step 2: insert a new row at row 2, pushing all existing data rows down ---
$ws.Rows.Item(2).Insert()

# Copy the formatting (styles/number formats) from the row immediately below
# (which now holds what used to be the old row 2) onto the newly inserted row.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)

# --- Step 3: populate the new row 2 with the latest day's data ---
# Force column A and E to plain text first so Excel does not reinterpret the
# dd-mm-yyyy strings as dates, then restore the original cell formatting
# (this keeps the same style index as the surrounding rows).
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "13-11-2025"
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Cells.Item(2, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(2, 3).Value = "IE07"
$ws.Cells.Item(2, 4).Value = 297.15

$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "01-11-2025"
$ws.Range("E3").Copy()
$ws.Range("E2").PasteSpecial(-4122)

# --- Step 4: rebuild the Circular Link hyperlinks for column F. ---
# Inserting the row above does not shift the worksheet's hyperlink
# relationships along with the data, so clear out whatever remains and
# recreate the complete set of hyperlinks in the correct, now-shifted order.
$ws.Range("F2").Hyperlinks.Delete()

$newUrl = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-11-2025.pdf"
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), $newUrl)

$destRow = 3
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Hyperlinks.Add($ws.Cells.Item($destRow, 6), $urls[$r])
    $destRow = $destRow + 1
}

# Adding hyperlinks registers a built-in "Hyperlink" cell style (blue,
# underlined) and re-styles every touched cell with it. The source data
# keeps the plain table styling instead, so drop the named style again and
# restore column F's formatting. Column A shares the exact same style as
# column F and was never touched by the hyperlink calls above, so it is a
# safe, still-correctly-styled source for the PasteSpecial.
$wb.Styles.Item("Hyperlink").Delete()
$ws.Range("A3").Copy()
$ws.Range("F2:F100").PasteSpecial(-4122)
